# Insert a new "Summary" slide just before the final "Thank You!" slide.
#
# Before:  ... , Slide21, Slide22 ("Thank You!")
# After:   ... , Slide21, Slide22 ("Summary"), Slide23 ("Thank You!")
#
# The existing "Thank You!" slide (SlideID 277) keeps its id/content and is
# simply pushed one position later; the brand-new slide gets the next free
# SlideID (278) and is inserted at index 22, matching the target deck's
# <p:sldIdLst> ordering (…278, 277 at the end).

$p = $ppt.ActivePresentation

# Index 22 = just before the current last slide ("Thank You!"); layout 2 =
# ppLayoutText ("Title and Content"), matching the slideLayout2.xml layout
# used by the new slide (title placeholder + idx=1 body placeholder).
$s = $p.Slides.Add(22, 2)

# --- Title placeholder -----------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Summary"

# --- Body / content placeholder ---------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange

$para1 = "So after doing some basic exploratory data analysis and Feature Engineering along with Feature Selection then finally clustering we were able to segment places and provide insights"
$para2 = "This project can help various people if they intent to "
$para3 = "Start a new business and want to get some idea"
$para4 = "An established business trying to find out their competition"
$para5 = "Probability of success of business"
$para6 = "Missing business in the locality "

# Build the paragraphs one at a time (instead of one multi-line Text
# assignment) so every run gets its own <a:rPr>.
$body.Text = $para1
[void]$body.InsertAfter("`r" + $para2)
[void]$body.InsertAfter("`r" + $para3)
[void]$body.InsertAfter("`r" + $para4)
[void]$body.InsertAfter("`r" + $para5)
[void]$body.InsertAfter("`r" + $para6)
# Trailing empty, numbered paragraph at the end of the list.
[void]$body.InsertAfter("`r")

# Paragraphs 3-7 (the "Start a new business" ... through the trailing blank
# paragraph) form an arabic-numbered list using the major-latin theme font
# for the bullet glyph, matching the target's <a:buFont typeface="+mj-lt"/>
# <a:buAutoNum type="arabicPeriod"/>.
for ($i = 3; $i -le 7; $i++) {
    $para = $body.Paragraphs($i)
    $para.ParagraphFormat.Bullet.Font.Name = "+mj-lt"
    $para.ParagraphFormat.Bullet.Type = 2
    $para.ParagraphFormat.Bullet.Style = 3
}
